$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing "field_wbddh_data_type" block (row 28)
# so that a new "Other" option is added as the first entry, shifting the
# existing field_wbddh_data_type / field_granularity_list rows down by one.
$ws.Rows("28:28").Insert()

# New row 28: field_wbddh_data_type / Other / Other
$ws.Range("A28").Value = "field_wbddh_data_type"
$ws.Range("B28").Value = "Other"
$ws.Range("C28").Value = "Other"

# The list_value_name (column C) for every field_wbddh_data_type row is now
# a constant "Other" value.
$ws.Range("C29").Value = "Other"
$ws.Range("C30").Value = "Other"
$ws.Range("C31").Value = "Other"
$ws.Range("C32").Value = "Other"

# Reflect the new active selection left behind by the edit.
$ws.Activate() | Out-Null
$ws.Range("B29").Select() | Out-Null
